$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value2 = 'Dialogue Name'
$ws.Range("B1").Value2 = 'Character'
$ws.Range("C1").Value2 = 'Line'
$ws.Range("D1").Value2 = 'Priority'
$ws.Range("E1").Value2 = 'Repeatable'
$ws.Range("F1").Value2 = 'Conditions'
$ws.Range("G1").Value2 = 'Remember'
$ws.Range("H1").Value2 = 'Context'

# --- Row 2 ---
$ws.Range("B2").Value2 = 'None'
$ws.Range("C2").Value2 = 'AtigMDEYpZVFApUVAPNqciSDkTYzgLZyyzbgQEbgMjwLXuwadAcxckUCTmE.'
$ws.Range("D2").Value2 = 0
$ws.Range("E2").Value2 = '-'
$ws.Range("F2").Value2 = 'None'
$ws.Range("G2").Value2 = 'None'
$ws.Range("H2").Value2 = 'Nonsense to get 120 character string'

# --- Row 3 ---
$ws.Range("A3").Value2 = 'Observed_Sign'
$ws.Range("B3").Value2 = 'Detective'
$ws.Range("C3").Value2 = '"Big Time Crime Park"...'
$ws.Range("D3").Value2 = 0
$ws.Range("E3").Value2 = 'Yes'
$ws.Range("F3").Value2 = 'None'
$ws.Range("G3").Value2 = 'None'
$ws.Range("H3").Value2 = 'Clicked on park sign'

# --- Row 4 ---
$ws.Range("A4").Value2 = 'Observed_Sign'
$ws.Range("B4").Value2 = 'Detective'
$ws.Range("C4").Value2 = 'Been getting a lot of reports about unusual happenings here.'
$ws.Range("D4").Value2 = 0
$ws.Range("E4").Value2 = '-'
$ws.Range("F4").Value2 = 'None'
$ws.Range("G4").Value2 = 'None'
$ws.Range("H4").Value2 = 'prev'

# --- Row 5 ---
$ws.Range("A5").Value2 = 'Hi_Charlie'
$ws.Range("B5").Value2 = 'Detective'
$ws.Range("C5").Value2 = 'Excuse me, have you seen any aliens around?'
$ws.Range("D5").Value2 = 1
$ws.Range("E5").Value2 = 'No'
$ws.Range("F5").Value2 = 'None'
$ws.Range("G5").Value2 = 'None'
$ws.Range("H5").Value2 = 'Clicked on chalk outline'

# --- Row 6 ---
$ws.Range("A6").Value2 = 'Hi_Charlie'
$ws.Range("B6").Value2 = 'Charlie'
$ws.Range("C6").Value2 = 'What? Uh, no...'
$ws.Range("D6").Value2 = 0
$ws.Range("E6").Value2 = '-'
$ws.Range("F6").Value2 = 'None'
$ws.Range("G6").Value2 = 'None'
$ws.Range("H6").Value2 = 'prev'

# --- Row 7 ---
$ws.Range("A7").Value2 = 'Hi_Charlie'
$ws.Range("B7").Value2 = 'Detective'
$ws.Range("C7").Value2 = 'What happened to you then?'
$ws.Range("D7").Value2 = 0
$ws.Range("E7").Value2 = '-'
$ws.Range("F7").Value2 = 'None'
$ws.Range("G7").Value2 = 'None'
$ws.Range("H7").Value2 = 'prev'

# --- Row 8 ---
$ws.Range("A8").Value2 = 'Hi_Charlie'
$ws.Range("B8").Value2 = 'Charlie'
$ws.Range("C8").Value2 = 'I... Uh... Fell from a high place. Completely by accident!'
$ws.Range("D8").Value2 = 0
$ws.Range("E8").Value2 = '-'
$ws.Range("F8").Value2 = 'None'
$ws.Range("G8").Value2 = 'None'
$ws.Range("H8").Value2 = 'prev'

# --- Row 9 ---
$ws.Range("A9").Value2 = 'Hi_Charlie'
$ws.Range("B9").Value2 = 'Detective'
$ws.Range("C9").Value2 = 'I see.'
$ws.Range("D9").Value2 = 0
$ws.Range("E9").Value2 = '-'
$ws.Range("F9").Value2 = 'None'
$ws.Range("G9").Value2 = 'None'
$ws.Range("H9").Value2 = 'prev'

# --- Row 10 ---
$ws.Range("A10").Value2 = 'Wait_Charlie'
$ws.Range("B10").Value2 = 'Detective'
$ws.Range("C10").Value2 = 'Could you stay here and keep a lookout for any aliens for me?'
$ws.Range("D10").Value2 = -1
$ws.Range("E10").Value2 = 'Yes'
$ws.Range("F10").Value2 = 'None'
$ws.Range("G10").Value2 = 'None'
$ws.Range("H10").Value2 = 'Clicked on chalk outline'

# --- Row 11 ---
$ws.Range("A11").Value2 = 'Wait_Charlie'
$ws.Range("B11").Value2 = 'Charlie'
$ws.Range("C11").Value2 = 'I, um, I''m not going anywhere.'
$ws.Range("D11").Value2 = 0
$ws.Range("E11").Value2 = '-'
$ws.Range("F11").Value2 = 'None'
$ws.Range("G11").Value2 = 'None'
$ws.Range("H11").Value2 = 'prev'

# --- Row 12 ---
$ws.Range("A12").Value2 = 'Wait_Charlie'
$ws.Range("B12").Value2 = 'Detective'
$ws.Range("C12").Value2 = 'Thanks!'
$ws.Range("D12").Value2 = 0
$ws.Range("E12").Value2 = '-'
$ws.Range("F12").Value2 = 'None'
$ws.Range("G12").Value2 = 'None'
$ws.Range("H12").Value2 = 'prev'

# --- Row 13 ---
$ws.Range("A13").Value2 = 'Rand_Hench_1'
$ws.Range("B13").Value2 = 'Henchman'
$ws.Range("C13").Value2 = 'Keep it moving, bub.'
$ws.Range("D13").Value2 = -1
$ws.Range("E13").Value2 = 'Yes'
$ws.Range("F13").Value2 = 'None'
$ws.Range("G13").Value2 = 'None'
$ws.Range("H13").Value2 = 'Clicked on henchman'

# --- Row 14 ---
$ws.Range("A14").Value2 = 'Rand_Hench_2'
$ws.Range("B14").Value2 = 'Detective'
$ws.Range("C14").Value2 = 'I have a few safety concerns about your park.'
$ws.Range("D14").Value2 = -1
$ws.Range("E14").Value2 = 'Yes'
$ws.Range("F14").Value2 = 'None'
$ws.Range("G14").Value2 = 'None'
$ws.Range("H14").Value2 = 'Clicked on henchman'

# --- Row 15 ---
$ws.Range("A15").Value2 = 'Rand_Hench_2'
$ws.Range("B15").Value2 = 'Henchman'
$ws.Range("C15").Value2 = 'Your concerns are noted.'
$ws.Range("D15").Value2 = 0
$ws.Range("E15").Value2 = '-'
$ws.Range("F15").Value2 = 'None'
$ws.Range("G15").Value2 = 'None'
$ws.Range("H15").Value2 = 'prev'

# --- Row 16 ---
$ws.Range("A16").Value2 = 'Hi_Bert'
$ws.Range("B16").Value2 = 'Detective'
$ws.Range("C16").Value2 = 'Excuse me, small child, have you seen any aliens around?'
$ws.Range("D16").Value2 = 0
$ws.Range("E16").Value2 = 'No'
$ws.Range("F16").Value2 = 'None'
$ws.Range("G16").Value2 = 'None'
$ws.Range("H16").Value2 = 'Clicked on Bert'

# --- Row 17 ---
$ws.Range("A17").Value2 = 'Hi_Bert'
$ws.Range("B17").Value2 = 'Bert'
$ws.Range("C17").Value2 = 'Haven''t seen one, bub. No aliens here. Nope. But...'
$ws.Range("D17").Value2 = 0
$ws.Range("E17").Value2 = '-'
$ws.Range("F17").Value2 = 'None'
$ws.Range("G17").Value2 = 'None'
$ws.Range("H17").Value2 = 'prev'

# --- Row 18 ---
$ws.Range("A18").Value2 = 'Hi_Bert'
$ws.Range("B18").Value2 = 'Bert'
$ws.Range("C18").Value2 = 'I did lose one of my "very special balloons". Green one.'
$ws.Range("D18").Value2 = 0
$ws.Range("E18").Value2 = '-'
$ws.Range("F18").Value2 = 'None'
$ws.Range("G18").Value2 = 'None'
$ws.Range("H18").Value2 = 'prev'

# --- Row 19 ---
$ws.Range("A19").Value2 = 'Hi_Bert'
$ws.Range("B19").Value2 = 'Bert'
$ws.Range("C19").Value2 = 'Find my balloon and I''ll make it worth your while.'
$ws.Range("D19").Value2 = 0
$ws.Range("E19").Value2 = '-'
$ws.Range("F19").Value2 = 'None'
$ws.Range("G19").Value2 = 'None'
$ws.Range("H19").Value2 = 'prev'

# --- Row 20 ---
$ws.Range("A20").Value2 = 'Bert_NoBalloon'
$ws.Range("B20").Value2 = 'Bert'
$ws.Range("C20").Value2 = 'Yo, you found my "special balloon" yet?'
$ws.Range("D20").Value2 = 0
$ws.Range("E20").Value2 = 'Yes'
$ws.Range("F20").Value2 = 'not has_balloon'
$ws.Range("G20").Value2 = 'None'
$ws.Range("H20").Value2 = 'Clicked on Bert'

# --- Row 21 ---
$ws.Range("A21").Value2 = 'Bert_NoBalloon'
$ws.Range("B21").Value2 = 'Detective'
$ws.Range("C21").Value2 = 'Not yet.'
$ws.Range("D21").Value2 = 0
$ws.Range("E21").Value2 = '-'
$ws.Range("F21").Value2 = 'None'
$ws.Range("G21").Value2 = 'None'
$ws.Range("H21").Value2 = 'prev'

# --- Row 22 ---
$ws.Range("A22").Value2 = 'Bert_HasBalloon'
$ws.Range("B22").Value2 = 'Bert'
$ws.Range("C22").Value2 = 'Yo, you found my "special balloon" yet?'
$ws.Range("D22").Value2 = 10
$ws.Range("E22").Value2 = 'No'
$ws.Range("F22").Value2 = 'has_balloon'
$ws.Range("G22").Value2 = 'None'
$ws.Range("H22").Value2 = 'Clicked on Bert'

# --- Row 23 ---
$ws.Range("A23").Value2 = 'Bert_HasBalloon'
$ws.Range("B23").Value2 = 'Detective'
$ws.Range("C23").Value2 = 'Yes! Here it is.'
$ws.Range("D23").Value2 = 0
$ws.Range("E23").Value2 = '-'
$ws.Range("F23").Value2 = 'None'
$ws.Range("G23").Value2 = 'None'
$ws.Range("H23").Value2 = 'prev'

# --- Row 24 ---
$ws.Range("A24").Value2 = 'Bert_HasBalloon'
$ws.Range("B24").Value2 = 'Bert'
$ws.Range("C24").Value2 = '[i]Looks both ways[/i]'
$ws.Range("D24").Value2 = 0
$ws.Range("E24").Value2 = '-'
$ws.Range("F24").Value2 = 'None'
$ws.Range("G24").Value2 = 'None'
$ws.Range("H24").Value2 = 'prev'
$ws.Rows.Item(24).RowHeight = 30.5

# --- Row 25 ---
$ws.Range("A25").Value2 = 'Bert_HasBalloon'
$ws.Range("B25").Value2 = 'Bert'
$ws.Range("C25").Value2 = '[i]Slips you a twenty[/i]'
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = '-'
$ws.Range("F25").Value2 = 'None'
$ws.Range("G25").Value2 = 'None'
$ws.Range("H25").Value2 = 'prev'
$ws.Rows.Item(25).RowHeight = 30.5

# --- Row 26 ---
$ws.Range("A26").Value2 = 'Bert_HasBalloon'
$ws.Range("B26").Value2 = 'Bert'
$ws.Range("C26").Value2 = 'You never saw me. Got it?'
$ws.Range("D26").Value2 = 0
$ws.Range("E26").Value2 = '-'
$ws.Range("F26").Value2 = 'None'
$ws.Range("G26").Value2 = 'None'
$ws.Range("H26").Value2 = 'prev'
$ws.Rows.Item(26).RowHeight = 30.5

# --- Row 27 ---
$ws.Range("A27").Value2 = 'Bert_WrongBalloon'
$ws.Range("B27").Value2 = 'Bert'
$ws.Range("C27").Value2 = 'What is that? A red herring? No! I lost my green balloon!'
$ws.Range("D27").Value2 = 1
$ws.Range("E27").Value2 = 'No'
$ws.Range("F27").Value2 = 'has_redBalloon'
$ws.Range("G27").Value2 = 'None'
$ws.Range("H27").Value2 = 'Clicked on Bert'
$ws.Rows.Item(27).RowHeight = 30.5

# --- Row 28 ---
$ws.Range("A28").Value2 = 'Sam_Hi'
$ws.Range("B28").Value2 = 'Sam'
$ws.Range("C28").Value2 = 'Uuuurgh...'
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 'No'
$ws.Range("F28").Value2 = 'None'
$ws.Range("G28").Value2 = 'None'
$ws.Range("H28").Value2 = 'Clicked on Sam'
$ws.Rows.Item(28).RowHeight = 30.5

# --- Row 29 ---
$ws.Range("A29").Value2 = 'Sam_Hi'
$ws.Range("B29").Value2 = 'Detective'
$ws.Range("C29").Value2 = 'You''re not looking so good...'
$ws.Range("D29").Value2 = 0
$ws.Range("E29").Value2 = '-'
$ws.Range("F29").Value2 = 'None'
$ws.Range("G29").Value2 = 'None'
$ws.Range("H29").Value2 = 'prev'
$ws.Rows.Item(29).RowHeight = 30.5

# --- Row 30 ---
$ws.Range("A30").Value2 = 'Sam_Hi'
$ws.Range("B30").Value2 = 'Detective'
$ws.Range("C30").Value2 = 'Did the aliens get you?'
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = '-'
$ws.Range("F30").Value2 = 'None'
$ws.Range("G30").Value2 = 'None'
$ws.Range("H30").Value2 = 'prev'
$ws.Rows.Item(30).RowHeight = 30.5

# --- Row 31 ---
$ws.Range("A31").Value2 = 'Sam_Hi'
$ws.Range("B31").Value2 = 'Sam'
$ws.Range("C31").Value2 = 'What''s it.. To you?'
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = '-'
$ws.Range("F31").Value2 = 'None'
$ws.Range("G31").Value2 = 'None'
$ws.Range("H31").Value2 = 'prev'
$ws.Rows.Item(31).RowHeight = 30.5

# --- Row 32 ---
$ws.Range("A32").Value2 = 'Sam_Hi'
$ws.Range("B32").Value2 = 'Sam'
$ws.Range("C32").Value2 = 'Listen.. get me some human medicine and I''ll tell you.'
$ws.Range("D32").Value2 = 0
$ws.Range("E32").Value2 = '-'
$ws.Range("F32").Value2 = 'None'
$ws.Range("G32").Value2 = 'None'
$ws.Range("H32").Value2 = 'prev'
$ws.Rows.Item(32).RowHeight = 30.5

# --- Column width (A), view/pane/selection updates ---
$ws.Columns.Item(1).ColumnWidth = 20

$win = $ws.Application.ActiveWindow
$win.ScrollRow = 8
$win.ScrollColumn = 1
[void]$ws.Range("C12").Select()

